# 自动更新Excel文件 - 2026-02-22 23:17:48
#
# Weekly rollover of the water-delivery tracking sheet: "today" advances
# from 2026-02-22 to 2026-02-23. For every data row (2..99), column E
# ("剩余" = days remaining) counts down by 1. Once it would hit 0 the
# cycle restarts: E is reset to the row's total (column D, "总天") and
# the start date in column F ("开始时间") is reset to the new "today"
# (20260223). Row 36 is a known bad/garbled record (its F value isn't a
# real date and E already equals D) so it is left untouched, matching
# the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newToday = 20260223
$firstRow = 2
$lastRow = 99
$skipRow = 36

for ($row = $firstRow; $row -le $lastRow; $row++) {
    if ($row -eq $skipRow) {
        continue
    }

    $totalDays = $ws.Cells.Item($row, 4).Value()
    $remaining = $ws.Cells.Item($row, 5).Value()

    $newRemaining = $remaining - 1

    if ($newRemaining -le 0) {
        $ws.Cells.Item($row, 5).Value = $totalDays
        $ws.Cells.Item($row, 6).Value = $newToday
    } else {
        $ws.Cells.Item($row, 5).Value = $newRemaining
    }
}
